$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "back to initial version of topics": column B ("Label") had been
# edited to an underscore_cased variant that lived in a separate set of
# shared-string entries; restore it to the same spaced text already used
# in column C ("Alpha") for those rows, so the duplicate shared strings
# collapse back down again.
$rows = @(6,7,12,13,14,25,28,29,33,38,40,41,42,45,46,47,52,53,54,58,59,67,68,75,76,78,79,80)
foreach ($r in $rows) {
  $c = $ws.Range("C$r").Value2
  $ws.Range("B$r").Value = $c
}

# Restore the view state captured in the saved workbook: scrolled down so
# row 17 is at the top, with the whole of column B selected (active cell
# sitting on B26).
[void]$ws.Range("B1:B1048576").Select()
